# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
# Updates Price (D) and Volume(1h) (E) columns, plus a row swap for Bittensor/dogwifhat (rows 40-41)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.497.29"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.002.99"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "2.993.52"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +3.03%  "
$ws.Range("E11").Value = "  +7.13%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D17").Value = "3.501.86"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "3.001.41"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "59.434.99"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +11.51%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +5.87%  "
$ws.Range("E35").Value = "  +5.81%  "
$ws.Range("D36").Value = "0.0₃0764"
$ws.Range("E36").Value = "  +9.01%  "
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E40").Value = "  +6.49%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E41").Value = "  +8.52%  "
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "2.780.61"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +21.17%  "
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  +0.16%  "

# The following Price cells contain values that look numeric (e.g. "563.15").
# Excel would auto-convert a plain string assignment into a real number, which
# would change the cell type away from text. Force a text number format while
# assigning, then clear the transient formatting so the cell keeps its original
# (unstyled) appearance but remains text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.15"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.38"
$ws.Range("D6").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.84"
$ws.Range("D14").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.36"
$ws.Range("D16").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.72"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.67"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.724"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.53"
$ws.Range("D24").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("D27").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.83"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.100"
$ws.Range("D33").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.02"
$ws.Range("D38").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "411.77"
$ws.Range("D41").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.00"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.68"
$ws.Range("D48").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.64"
$ws.Range("D51").ClearFormats()
